$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation for 2026/01/21 (11:00, ranking 21) needs to be inserted as
# row 668, pushing the existing rows 668-709 down to 669-710.
#
# Row 667 already holds a 2026/01/21 entry, so we copy it (to inherit its exact
# text formatting for the date/weekday columns) down into a freshly inserted
# row 668, then overwrite the time/ranking columns (C/D) with the new values.
$ws.Rows.Item(667).Copy()
$ws.Rows.Item(668).Insert()

$ws.Cells.Item(668, 3).Value = 11
$ws.Cells.Item(668, 4).Value = 21
